# Automatic update of files.
# Applies the targeted cell edits described by the upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple "Taxonsorteringsordning" (column B) bumps on otherwise
#     untouched rows -------------------------------------------------
$ws.Range("B2").Value  = 79245
$ws.Range("B3").Value  = 79245
$ws.Range("B7").Value  = 79245
$ws.Range("B9").Value  = 79245
$ws.Range("B10").Value = 79245
$ws.Range("B11").Value = 79245
$ws.Range("B12").Value = 79245
$ws.Range("B13").Value = 79245
$ws.Range("B14").Value = 79245
$ws.Range("B15").Value = 79245
$ws.Range("B16").Value = 91831
$ws.Range("B17").Value = 79245

# --- Rows 4 and 5 swap their record content (species data), each also
#     getting its own distinct new Taxonsorteringsordning value, and the
#     "Publik kommentar" (AC) note moves from row 4 to row 5 -----------

# Row 4 becomes what used to be row 5's species record.
$ws.Range("A4").Value = 131136874
$ws.Range("B4").Value = 79245
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("J4").Value = "bålar"
$ws.Range("Q4").Value = 788960
$ws.Range("R4").Value = 7131416
$ws.Range("AC4").ClearContents()

# Row 5 becomes what used to be row 4's species record.
$ws.Range("A5").Value = 131136941
$ws.Range("B5").Value = 83091
$ws.Range("E5").Value = 1312
$ws.Range("F5").Value = "Gammelgransskål"
$ws.Range("G5").Value = "Pseudographis pinicola"
$ws.Range("H5").Value = "(Nyl.) Rehm"
$ws.Range("J5").Value = "fruktkroppar"
$ws.Range("Q5").Value = 788995
$ws.Range("R5").Value = 7131220
$ws.Range("AC5").Value = "på en gammal senvuxen gran"
